$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (rows 2 and 3),
# pushing every existing data row down by two positions.
$ws.Rows("2:3").Insert()

# The newly inserted rows inherit the header row's formatting by default;
# restore the same formatting as the rest of the data block (row 4, which
# used to be row 2) by copying formats down onto the new rows. Copy only
# the used A:T range (not the whole row) so the sheet's used range does
# not balloon out to column XFD.
$ws.Range("A4:T4").Copy()
$ws.Range("A2:T3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row 2 with the new weekly record.
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C2").Value = "Ñuble"
$ws.Range("D2").Value = 45257
$ws.Range("E2").Value = 16
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100101
$ws.Range("H2").Value = "Berries"
$ws.Range("I2").Value = 100101001
$ws.Range("J2").Value = "Arándano (blue)"
$ws.Range("K2").Value = "Sin especificar"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 6000
$ws.Range("O2").Value = 6000
$ws.Range("P2").Value = 6000
$ws.Range("Q2").Value = '$/bandeja 2 kilos'
$ws.Range("R2").Value = "Región de Ñuble"
$ws.Range("S2").Value = 3000
$ws.Range("T2").Value = 2

# Populate the new row 3 with the second new weekly record.
$ws.Range("A3").Value = 7
$ws.Range("B3").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C3").Value = "Ñuble"
$ws.Range("D3").Value = 45257
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100101
$ws.Range("H3").Value = "Berries"
$ws.Range("I3").Value = 100101001
$ws.Range("J3").Value = "Arándano (blue)"
$ws.Range("K3").Value = "Sin especificar"
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 5000
$ws.Range("O3").Value = 5000
$ws.Range("P3").Value = 5000
$ws.Range("Q3").Value = '$/bandeja 2 kilos'
$ws.Range("R3").Value = "Región de Ñuble"
$ws.Range("S3").Value = 2500
$ws.Range("T3").Value = 2
